$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update probability matrix values per games pulled March 7

# Row 2
$ws.Range("B2").Value = 0.1384615384615385
$ws.Range("C2").Value = 0.6538461538461539
$ws.Range("J2").Value = 0.01538461538461539
$ws.Range("O2").Value = 0.003846153846153846
$ws.Range("P2").Value = 0.1076923076923077
$ws.Range("S2").Value = 0.08076923076923077

# Row 3
$ws.Range("C3").Value = 0.0594059405940594
$ws.Range("J3").Value = 0.0198019801980198
$ws.Range("P3").Value = 0.801980198019802
$ws.Range("S3").Value = 0.1188118811881188

# Row 4
$ws.Range("J4").Value = 0.1090909090909091
$ws.Range("P4").Value = 0.6727272727272727
$ws.Range("S4").Value = 0.2181818181818182

# Row 6
$ws.Range("B6").Value = 0.04020100502512563
$ws.Range("D6").Value = 0.02010050251256281
$ws.Range("E6").Value = 0.005025125628140704
$ws.Range("F6").Value = 0.05025125628140704
$ws.Range("J6").Value = 0.2512562814070352
$ws.Range("O6").Value = 0.01005025125628141
$ws.Range("Q6").Value = 0.221105527638191
$ws.Range("R6").Value = 0.08542713567839195
$ws.Range("S6").Value = 0.3165829145728643

# Row 7
$ws.Range("B7").Value = 0.09947643979057591
$ws.Range("D7").Value = 0.01570680628272251
$ws.Range("F7").Value = 0.05235602094240838
$ws.Range("J7").Value = 0.1465968586387434
$ws.Range("O7").Value = 0.01570680628272251
$ws.Range("Q7").Value = 0.1727748691099476
$ws.Range("R7").Value = 0.06806282722513089
$ws.Range("S7").Value = 0.4293193717277487

# Row 8
$ws.Range("B8").Value = 0.0851063829787234
$ws.Range("D8").Value = 0.02836879432624113
$ws.Range("E8").Value = 0.002364066193853428
$ws.Range("F8").Value = 0.05673758865248227
$ws.Range("J8").Value = 0.09692671394799054
$ws.Range("O8").Value = 0.01891252955082742
$ws.Range("Q8").Value = 0.2033096926713948
$ws.Range("R8").Value = 0.132387706855792
$ws.Range("S8").Value = 0.375886524822695

# Row 9
$ws.Range("B9").Value = 0.07727272727272727
$ws.Range("D9").Value = 0.00909090909090909
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.1136363636363636
$ws.Range("O9").Value = 0.01818181818181818
$ws.Range("Q9").Value = 0.2318181818181818
$ws.Range("R9").Value = 0.1045454545454545
$ws.Range("S9").Value = 0.3954545454545454

# Row 10
$ws.Range("B10").Value = 0.104570100697134
$ws.Range("D10").Value = 0.02478698683191325
$ws.Range("F10").Value = 0.05886909372579396
$ws.Range("J10").Value = 0.1247095274980635
$ws.Range("O10").Value = 0.0139426800929512
$ws.Range("Q10").Value = 0.2215336948102246
$ws.Range("R10").Value = 0.09450038729666925
$ws.Range("S10").Value = 0.3570875290472502

# Row 11
$ws.Range("G11").Value = 0.1360544217687075
$ws.Range("J11").Value = 0.09863945578231292
$ws.Range("K11").Value = 0.1768707482993197
$ws.Range("L11").Value = 0.5816326530612245
$ws.Range("S11").Value = 0.006802721088435374

# Row 12
$ws.Range("G12").Value = 0.7267441860465116
$ws.Range("J12").Value = 0.2325581395348837
$ws.Range("L12").Value = 0.01744186046511628
$ws.Range("S12").Value = 0.02325581395348837

# Row 13
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.25

# Row 15
$ws.Range("F15").Value = 0.03448275862068965
$ws.Range("H15").Value = 0.1724137931034483
$ws.Range("I15").Value = 0.07389162561576355
$ws.Range("J15").Value = 0.3103448275862069
$ws.Range("K15").Value = 0.09359605911330049
$ws.Range("M15").Value = 0.004926108374384237
$ws.Range("O15").Value = 0.0541871921182266
$ws.Range("S15").Value = 0.2561576354679803

# Row 16
$ws.Range("F16").Value = 0.0045662100456621
$ws.Range("H16").Value = 0.1461187214611872
$ws.Range("I16").Value = 0.0639269406392694
$ws.Range("J16").Value = 0.4018264840182648
$ws.Range("K16").Value = 0.1598173515981735
$ws.Range("M16").Value = 0.0273972602739726
$ws.Range("N16").Value = 0.0045662100456621
$ws.Range("O16").Value = 0.0547945205479452
$ws.Range("S16").Value = 0.136986301369863

# Row 17
$ws.Range("F17").Value = 0.01364522417153996
$ws.Range("H17").Value = 0.189083820662768
$ws.Range("I17").Value = 0.0935672514619883
$ws.Range("J17").Value = 0.4113060428849902
$ws.Range("K17").Value = 0.1072124756335283
$ws.Range("M17").Value = 0.01364522417153996
$ws.Range("O17").Value = 0.05847953216374269
$ws.Range("S17").Value = 0.1130604288499025

# Row 18
$ws.Range("F18").Value = 0.01716738197424893
$ws.Range("H18").Value = 0.1630901287553648
$ws.Range("I18").Value = 0.07296137339055794
$ws.Range("J18").Value = 0.4978540772532189
$ws.Range("K18").Value = 0.09012875536480687
$ws.Range("M18").Value = 0.0128755364806867
$ws.Range("O18").Value = 0.04291845493562232
$ws.Range("S18").Value = 0.1030042918454936

# Row 19
$ws.Range("F19").Value = 0.01812191103789127
$ws.Range("H19").Value = 0.1836902800658979
$ws.Range("I19").Value = 0.1046128500823723
$ws.Range("J19").Value = 0.3723228995057661
$ws.Range("K19").Value = 0.09390444810543658
$ws.Range("M19").Value = 0.01729818780889621
$ws.Range("N19").Value = 0.0008237232289950577
$ws.Range("O19").Value = 0.06836902800658978
$ws.Range("S19").Value = 0.1408566721581549
